# Update CCB data using new condition list
#
# The physical worksheet "sheet1.xml" in the package corresponds to the
# sheet named "variableNames" (sheetId 1 / rId1) - NOT the sheet literally
# named "Sheet1" (which is sheetId 3 / rId3). All of the edits described
# by the diff live in that "variableNames" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variableNames")

# Insert a brand-new row above the current row 5 ("occupation"), shifting
# the old rows 5-28 down to 6-29, then populate the new row with the new
# "birthCountry" variable (field F29).
$ws.Rows("5:5").Insert()

$ws.Range("A5").Value = "birthCountry"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "F29"

# Match the new selection left behind in the sheet (E5) and make sure the
# "variableNames" tab stays the active one, as in the source file.
$ws.Activate()
$ws.Range("E5").Select()
